# df_scores.xlsx edit: rename Sheet1 -> all_scores, update its score arrays,
# and add a new "mean_scores" sheet summarizing the results.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the existing sheet and refresh its data ---------------------
$allScores = $wb.Worksheets.Item(1)
$allScores.Name = "all_scores"

$allScores.Range("B2").Value = "[0.67 0.84 0.67 0.71 0.86 0.73 0.69 0.69 0.8  0.7 ]"
$allScores.Range("C2").Value = "[0.69 0.85 0.65 0.58 0.76 0.6  0.56 0.6  0.8  0.76]"
$allScores.Range("D2").Value = "[0.64 0.84 0.68 0.84 0.96 0.85 0.81 0.77 0.8  0.64]"
$allScores.Range("E2").Value = "[0.75 0.89 0.77 0.74 0.93 0.77 0.79 0.74 0.9  0.82]"
$allScores.Range("F2").Value = "[0.51 0.78 0.54 0.48 0.86 0.53 0.58 0.48 0.79 0.64]"

$allScores.Range("B3").Value = "[0.69 0.78 0.69 0.76 0.84 0.63 0.69 0.69 0.82 0.78]"
$allScores.Range("C3").Value = "[0.73 0.77 0.65 0.69 0.76 0.52 0.6  0.48 0.8  0.76]"
$allScores.Range("D3").Value = "[0.64 0.8  0.72 0.84 0.92 0.73 0.77 0.88 0.84 0.8 ]"
$allScores.Range("E3").Value = "[0.77 0.84 0.74 0.84 0.95 0.69 0.77 0.77 0.88 0.84]"
$allScores.Range("F3").Value = "[0.54 0.68 0.48 0.68 0.9  0.38 0.54 0.54 0.75 0.69]"

# --- 2. Add the new "mean_scores" sheet right after "all_scores" -----------
$meanScores = $wb.Worksheets.Add($null, $allScores)
$meanScores.Name = "mean_scores"

# Headers (row 1) - copy same header text and formatting as all_scores
$allScores.Range("B1:F1").Copy()
$meanScores.Range("B1").PasteSpecial(-4163)
$meanScores.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row labels (column A) - copy formatting/text from all_scores
$allScores.Range("A2:A3").Copy()
$meanScores.Range("A2").PasteSpecial(-4163)
$meanScores.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows
$meanScores.Range("B2").Value = "0.73 +/- 0.07"
$meanScores.Range("C2").Value = "0.68 +/- 0.1"
$meanScores.Range("D2").Value = "0.78 +/- 0.1"
$meanScores.Range("E2").Value = "0.81 +/- 0.07"
$meanScores.Range("F2").Value = "0.62 +/- 0.14"

$meanScores.Range("B3").Value = "0.74 +/- 0.07"
$meanScores.Range("C3").Value = "0.68 +/- 0.11"
$meanScores.Range("D3").Value = "0.79 +/- 0.08"
$meanScores.Range("E3").Value = "0.81 +/- 0.07"
$meanScores.Range("F3").Value = "0.62 +/- 0.14"

$allScores.Activate()
